# Gonzales-08122024-08162024.docx edit
#
# Commit: "Add files via upload"
#   - Removes the run of text "IT Support Intern" that was sitting in the
#     "Department Assigned:" value cell of the weekly-report table, leaving
#     the paragraph (and its paragraph-mark formatting) empty.
#   - Drops the stray "_GoBack" bookmark left over from Word's last cursor
#     position tracking.
#
# (The rest of the upstream diff -- new wp14:anchorId/editId GUIDs on
# drawings, extra w16* namespace declarations, and one extra latent style
# entry -- are artifacts automatically stamped by the Word version that
# re-saved the file, not user edits, so there is nothing to replay for
# them here.)

$d = $word.ActiveDocument

# 1) Remove the "IT Support Intern" run from the Department Assigned cell.
#    Replacing with an empty string removes the <w:r> (and its <w:t>)
#    entirely, leaving the paragraph's <w:pPr> untouched -- exactly what
#    the diff shows.
$d.Content.Find.Execute(
    "IT Support Intern",   # FindText
    $true,                 # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                 # Forward
    1,                     # Wrap (wdFindContinue)
    $false,                # Format
    "",                    # ReplaceWith
    2                      # Replace (wdReplaceAll)
) | Out-Null

# 2) Best-effort cleanup of the orphaned "_GoBack" bookmark (a Word
#    "last edit location" marker with no visible content). Guarded so the
#    script still succeeds if the host doesn't surface it for deletion.
try {
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }
} catch {
    # Nothing addressable -- ignore.
}
